$d = $word.ActiveDocument

# The last paragraph in the document body is the Chinese translation
# paragraph ("我一直梦想...这门课。"). Append a new, empty paragraph
# right after it. InsertParagraphAfter() duplicates the paragraph/run
# formatting (spacing after=0, firstLine indent=420, justified, sz=21)
# of the paragraph it is called on, giving us the new blank paragraph
# with matching formatting and no text.
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$last.Range.InsertParagraphAfter()
